# Opus-Base Hyperparameter Search Results — add final opus-base fine adapt wce
# and some lsp simple results.
#
# Target sheets (workbook tab order):
#   1 opus_base Validation
#   2 opus_base Simple aWCE        (selection I26 -> J24)
#   3 opus_base AoN aWCE           (selection I23 -> J22)
#   4 opus_base Fine aWCE          (selection I15 -> I30, fills D/E/G, tab stays active)
#   5 opus_base LSP Simple aWCE    (selection I23 -> H20, fills C/D/F)
#   6 opus_base LSP AoN aWCE
#   7 opus_base LSP Fine aWCE

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "opus_base Fine aWCE": fill in Duration (D), WCE (E) and the final
# literal metric (G) for every row that was still blank. F holds a shared
# formula (=E/3600) so it recomputes automatically once E is populated.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("opus_base Fine aWCE")

$fineData = @(
    @{r=4; D=42.5869; E=19109.868900000001; G=103.4885},
    @{r=6; D=42.600099999999998; E=17395.145799999998; G=82.794399999999996},
    @{r=7; D=42.297699999999999; E=16430.728299999999; G=103.4885},
    @{r=8; D=42.670999999999999; E=17620.359799999998; G=82.794399999999996},
    @{r=9; D=42.6693; E=22515.164000000001; G=137.96010000000001},
    @{r=10; D=42.499400000000001; E=17713.3151; G=117.2388},
    @{r=11; D=42.659300000000002; E=20801.682199999999; G=103.4885},
    @{r=12; D=42.381799999999998; E=17628.0972; G=82.794399999999996},
    @{r=13; D=42.747199999999999; E=17665.2372; G=82.794399999999996},
    @{r=14; D=42.772199999999998; E=22714.4349; G=137.96010000000001},
    @{r=15; D=42.045999999999999; E=10507.374400000001; G=69.000799999999998},
    @{r=17; D=42.594099999999997; E=20454.816200000001; G=137.96010000000001},
    @{r=18; D=42.690899999999999; E=20444.894700000001; G=137.96010000000001},
    @{r=19; D=42.612000000000002; E=18350.4892; G=117.2388},
    @{r=20; D=42.553800000000003; E=14190.751399999999; G=103.4885},
    @{r=21; D=42.829099999999997; E=22787.285; G=137.96010000000001},
    @{r=22; D=42.479799999999997; E=15149.1173; G=103.4885},
    @{r=23; D=42.753100000000003; E=16079.589599999999; G=82.794399999999996},
    @{r=24; D=42.702199999999998; E=20685.9015; G=137.96010000000001},
    @{r=25; D=42.891399999999997; E=25843.4002; G=137.96010000000001},
    @{r=26; D=42.534399999999998; E=16114.2562; G=82.794399999999996},
    @{r=30; D=42.553899999999999; E=16043.049000000001; G=103.4885},
    @{r=33; D=42.909500000000001; E=25622.4293; G=137.96010000000001}
)

foreach ($row in $fineData) {
    $ws4.Range("D" + $row.r).Value = $row.D
    $ws4.Range("E" + $row.r).Value = $row.E
    $ws4.Range("G" + $row.r).Value = $row.G
}

# Rows 14 and 33 are the best ("final") result in their hyperparameter block,
# so the author highlighted the Duration cell with the green fill used
# elsewhere in the workbook for the winning configuration.
$ws4.Range("D14").Interior.Color = 5296274
$ws4.Range("D33").Interior.Color = 5296274

# ---------------------------------------------------------------------------
# Sheet "opus_base LSP Simple aWCE ": fill in Duration (C), WCE (D) and the
# final literal metric (F). E holds a shared formula (=D/3600).
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("opus_base LSP Simple aWCE ")

$lspSimpleData = @(
    @{r=2; C=42.409100000000002; D=15410.118399999999; F=82.794399999999996},
    @{r=3; C=42.9313; D=23419.246299999999; F=137.96010000000001},
    @{r=5; C=42.810600000000001; D=25720.1014; F=137.96010000000001},
    @{r=6; C=42.537199999999999; D=15398.700500000001; F=103.4885},
    @{r=7; C=42.713000000000001; D=17855.449400000001; F=82.794399999999996},
    @{r=8; C=42.569499999999998; D=17683.5304; F=82.794399999999996},
    @{r=10; C=42.671300000000002; D=18727.5838; F=117.2388},
    @{r=11; C=42.667700000000004; D=14459.339; F=103.4885},
    @{r=13; C=42.542999999999999; D=20668.981299999999; F=103.4885},
    @{r=14; C=42.842300000000002; D=25939.740300000001; F=137.96010000000001},
    @{r=15; C=42.505800000000001; D=15643.267; F=69.000799999999998},
    @{r=16; C=42.415199999999999; D=14178.332700000001; F=103.4885},
    @{r=17; C=42.814799999999998; D=25285.619900000002; F=137.96010000000001}
)

foreach ($row in $lspSimpleData) {
    $ws5.Range("C" + $row.r).Value = $row.C
    $ws5.Range("D" + $row.r).Value = $row.D
    $ws5.Range("F" + $row.r).Value = $row.F
}

# Rows 7 and 8 lose their top border (matching the border used one row up in
# this block) once they are filled in, same as the source workbook.
$ws5.Range("F7").Borders.Item(8).LineStyle = -4142
$ws5.Range("F8").Borders.Item(8).LineStyle = -4142

# ---------------------------------------------------------------------------
# Restore each sheet's last-used selection (activeCell) as left by the author.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("opus_base Simple aWCE")
$ws2.Activate()
$ws2.Range("J24").Select()

$ws3 = $wb.Worksheets.Item("opus_base AoN aWCE")
$ws3.Activate()
$ws3.Range("J22").Select()

$ws5.Activate()
$ws5.Range("H20").Select()

$ws4.Activate()
$ws4.Range("I30").Select()
